$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16, pushing the existing rows 16-41 down to 17-42.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row with the "Adobe Creative Cloud Experience"
# install-path entry (mirrors the pattern used by the surrounding rows).
$ws.Range("A16").Value = "%ProgramFilesFolder32%\Adobe\Adobe Creative Cloud Experience"
$ws.Range("E16").Value = "0x00000221"

# Move the active selection to A4, matching the saved view state.
$ws.Range("A4").Select()
